# InhalerWorking.pptx maintenance edit
#
# The "Preprocess / Process / PostProcess" pipeline diagram slide documents
# the old respiratory-circuit handling (Administer / RemoveActiveCommand /
# "Combined Circuit Solved by Respiratory" / "Combined Circuit Time Advanced
# by Respiratory"). That handling was reworked, so the now-stale diagram
# slide is removed from the deck as part of the documentation refresh.

$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isPipelineDiagramSlide = $false

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*Preprocess*") {
                $isPipelineDiagramSlide = $true
            }
        }
    }

    if ($isPipelineDiagramSlide) {
        $slide.Delete()
    }
}
